# Applies the cell-value changes described by the commit diff
# (regenerated data refresh for 北京-漫展信息.xlsx: 展览 + 全部类型 sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 2799
$ws.Cells.Item(5, 6).Value = 6375
$ws.Cells.Item(6, 6).Value = 2449
$ws.Cells.Item(9, 6).Value = 34
$ws.Cells.Item(10, 6).Value = 2878
$ws.Cells.Item(11, 6).Value = 38
$ws.Cells.Item(13, 6).Value = 6990
$ws.Cells.Item(19, 6).Value = 8113
$ws.Cells.Item(23, 6).Value = 58
$ws.Cells.Item(25, 6).Value = 18
$ws.Cells.Item(26, 6).Value = 115
$ws.Cells.Item(30, 6).Value = 34
$ws.Cells.Item(33, 6).Value = 2593
$ws.Cells.Item(39, 6).Value = 640
$ws.Cells.Item(40, 6).Value = 3649

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 2799
$ws.Cells.Item(7, 6).Value = 6375
$ws.Cells.Item(8, 6).Value = 2449
$ws.Cells.Item(12, 6).Value = 34
$ws.Cells.Item(13, 6).Value = 2878
$ws.Cells.Item(14, 6).Value = 38
$ws.Cells.Item(18, 6).Value = 6990
$ws.Cells.Item(24, 6).Value = 8113
$ws.Cells.Item(28, 6).Value = 58
$ws.Cells.Item(29, 3).Value = "北京·第16届IJOY漫展【樱桃专场见面会】"
$ws.Cells.Item(29, 6).Value = 18
$ws.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83461"
$ws.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/QhYUWCBC1711527705322.jpeg"
$ws.Cells.Item(30, 2).Value = "'2024-05-03"
$ws.Cells.Item(30, 3).Value = "北京·知名演员 川久保拓司 专场活动"
$ws.Cells.Item(30, 4).Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Cells.Item(30, 5).Value = "2024.05.03 10:30-05.03 15:00"
$ws.Cells.Item(30, 6).Value = 115
$ws.Cells.Item(30, 7).Value = 528
$ws.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82897"
$ws.Cells.Item(30, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/rxrJuuvX1710409029498.jpeg"
$ws.Cells.Item(31, 2).Value = "'2024-05-04"
$ws.Cells.Item(31, 3).Value = "北京·XW咒术回战only"
$ws.Cells.Item(31, 4).Value = "北花园路1号 超级蜂巢"
$ws.Cells.Item(31, 5).Value = "2024.05.04 10:00-05.04 17:00"
$ws.Cells.Item(31, 6).Value = 65
$ws.Cells.Item(31, 7).Value = 60
$ws.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83570"
$ws.Cells.Item(31, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/G9X2HmU11711703284044.jpeg"
$ws.Cells.Item(32, 3).Value = "北京·永劫无间only展"
$ws.Cells.Item(32, 4).Value = "中滩村北二街与立水桥东一路交叉口西南150米 天通苑街心花园"
$ws.Cells.Item(32, 5).Value = "2024.05.04 11:00-05.04 18:00"
$ws.Cells.Item(32, 6).Value = 32
$ws.Cells.Item(32, 7).Value = 70
$ws.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82350"
$ws.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/3fjr8Tll1709627193073.jpeg"
$ws.Cells.Item(34, 6).Value = 34
$ws.Cells.Item(38, 6).Value = 2593
$ws.Cells.Item(44, 6).Value = 640
$ws.Cells.Item(46, 6).Value = 3649
